$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $value) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

# Row 2
Set-TextValue 2 4 '97.674.56'

# Row 3
Set-TextValue 3 4 '3.726.51'
$ws.Cells.Item(3, 5).Value = '  +0.21%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.20%  '

# Row 5
Set-TextValue 5 4 '2.15'
$ws.Cells.Item(5, 5).Value = '  +10.69%  '

# Row 6
Set-TextValue 6 4 '237.26'
$ws.Cells.Item(6, 5).Value = '  -1.10%  '

# Row 7
Set-TextValue 7 4 '657.16'
$ws.Cells.Item(7, 5).Value = '  -0.11%  '

# Row 8
Set-TextValue 8 4 '0.445'
$ws.Cells.Item(8, 5).Value = '  +4.30%  '

# Row 9
Set-TextValue 9 4 '1.12'
$ws.Cells.Item(9, 5).Value = '  +2.65%  '

# Row 10
Set-TextValue 10 4 '1.00'
$ws.Cells.Item(10, 5).Value = '  +0.10%  '

# Row 11
Set-TextValue 11 4 '3.728.67'
$ws.Cells.Item(11, 5).Value = '  +0.20%  '

# Row 12
Set-TextValue 12 4 '0.0000314'
$ws.Cells.Item(12, 5).Value = '  +16.02%  '

# Row 13
Set-TextValue 13 4 '44.86'
$ws.Cells.Item(13, 5).Value = '  -1.65%  '

# Row 14
Set-TextValue 14 4 '0.207'

# Row 15
Set-TextValue 15 4 '6.84'
$ws.Cells.Item(15, 5).Value = '  -0.50%  '

# Row 16
Set-TextValue 16 4 '4.420.31'
$ws.Cells.Item(16, 5).Value = '  +0.25%  '

# Row 17
Set-TextValue 17 4 '97.433.33'
$ws.Cells.Item(17, 5).Value = '  +0.88%  '

# Row 18
Set-TextValue 18 4 '8.88'
$ws.Cells.Item(18, 5).Value = '  -3.16%  '

# Row 19
Set-TextValue 19 4 '3.712.60'
$ws.Cells.Item(19, 5).Value = '  +0.11%  '

# Row 20
Set-TextValue 20 4 '13.07'
$ws.Cells.Item(20, 5).Value = '  +0.97%  '

# Row 21
Set-TextValue 21 4 '18.98'
$ws.Cells.Item(21, 5).Value = '  -1.13%  '

# Row 22
Set-TextValue 22 4 '0.530'
$ws.Cells.Item(22, 5).Value = '  -0.67%  '

# Row 23
Set-TextValue 23 4 '529.34'
$ws.Cells.Item(23, 5).Value = '  +0.79%  '

# Row 24
Set-TextValue 24 4 '3.45'
$ws.Cells.Item(24, 5).Value = '  -1.60%  '

# Row 25
Set-TextValue 25 4 '0.0000225'
$ws.Cells.Item(25, 5).Value = '  +9.60%  '

# Row 26
Set-TextValue 26 4 '118.08'
$ws.Cells.Item(26, 5).Value = '  +14.58%  '

# Row 27
$ws.Cells.Item(27, 5).Value = '  -2.49%  '

# Row 28
Set-TextValue 28 4 '0.210'
$ws.Cells.Item(28, 5).Value = '  +24.58%  '

# Row 29
Set-TextValue 29 4 '13.49'
$ws.Cells.Item(29, 5).Value = '  +0.15%  '

# Row 30
Set-TextValue 30 4 '12.73'
$ws.Cells.Item(30, 5).Value = '  +0.10%  '

# Row 31
$ws.Cells.Item(31, 5).Value = '  -1.70%  '

# Row 32
Set-TextValue 32 4 '0.999'
$ws.Cells.Item(32, 5).Value = '  -0.11%  '

# Row 33
$ws.Cells.Item(33, 5).Value = '  +1.71%  '

# Row 34
$ws.Cells.Item(34, 2).Value = 'EthereumClassic'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue 34 4 '33.21'
$ws.Cells.Item(34, 5).Value = '  +0.43%  '

# Row 35
$ws.Cells.Item(35, 2).Value = 'Fetch.AI'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue 35 4 '1.82'
$ws.Cells.Item(35, 5).Value = '  -3.62%  '

# Row 36
Set-TextValue 36 4 '0.998'
$ws.Cells.Item(36, 5).Value = '  -0.36%  '

# Row 37
Set-TextValue 37 4 '0.598'
$ws.Cells.Item(37, 5).Value = '  -1.57%  '

# Row 38
Set-TextValue 38 4 '641.00'
$ws.Cells.Item(38, 5).Value = '  -2.85%  '

# Row 39
Set-TextValue 39 4 '8.77'
$ws.Cells.Item(39, 5).Value = '  -2.92%  '

# Row 40
$ws.Cells.Item(40, 5).Value = '  +0.02%  '

# Row 41
Set-TextValue 41 4 '0.168'
$ws.Cells.Item(41, 5).Value = '  +2.92%  '

# Row 42
$ws.Cells.Item(42, 2).Value = 'Filecoin'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 42 4 '6.84'
$ws.Cells.Item(42, 5).Value = '  -6.43%  '

# Row 43
$ws.Cells.Item(43, 2).Value = 'EnergySwap'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 43 4 '41.13'
$ws.Cells.Item(43, 5).Value = '  +2.14%  '

# Row 44
$ws.Cells.Item(44, 2).Value = 'Algorand'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue 44 4 '0.491'
$ws.Cells.Item(44, 5).Value = '  +8.71%  '

# Row 45
Set-TextValue 45 4 '2.01'
$ws.Cells.Item(45, 5).Value = '  -0.07%  '

# Row 46
Set-TextValue 46 4 '0.969'
$ws.Cells.Item(46, 5).Value = '  -1.39%  '

# Row 47
$ws.Cells.Item(47, 5).Value = '  -1.50%  '

# Row 48
Set-TextValue 48 4 '2.39'
$ws.Cells.Item(48, 5).Value = '  +1.29%  '

# Row 49
$ws.Cells.Item(49, 2).Value = 'WhiteBITCoin'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
Set-TextValue 49 4 '23.66'
$ws.Cells.Item(49, 5).Value = '  +0.10%  '

# Row 50
$ws.Cells.Item(50, 2).Value = 'Cosmos'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue 50 4 '8.76'
$ws.Cells.Item(50, 5).Value = '  +1.10%  '

# Row 51
Set-TextValue 51 4 '3.34'
$ws.Cells.Item(51, 5).Value = '  +2.90%  '
